# The document had three logo pictures whose auto-generated drawing
# "name" (wp:docPr / pic:cNvPr @name, e.g. "image1.png"/"image2.png")
# had drifted out of sync with the actual embedded media part names -
# a common artefact of copy/pasting the same logos between BTEC
# assignment-brief templates. This renumbers them back:
#
#   Footer primary (Pearson logo, id 3)  image2.png -> image1.png
#   Footer primary (Pearson logo, id 2)  image2.png -> image1.png
#   Header primary (BTEC logo,   id 1)  image1.jpg -> image2.jpg
#
# InlineShape has no writable .Name in the Word object model, so each
# picture is round-tripped through Shape (which does expose .Name) and
# back to an InlineShape, leaving its position/wrapping untouched.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineShape($inlineShape, [string]$newName) {
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape() | Out-Null
}

# --- Footers: both hold the Pearson Edexcel logo -----------------------
$footerA = $sec.Footers.Item(1)
Rename-InlineShape $footerA.Range.InlineShapes.Item(1) "image1.png"

$footerB = $sec.Footers.Item(2)
Rename-InlineShape $footerB.Range.InlineShapes.Item(1) "image1.png"

# --- Header: holds the BTEC logo ---------------------------------------
$header = $sec.Headers.Item(2)
Rename-InlineShape $header.Range.InlineShapes.Item(1) "image2.jpg"
